$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.3236798569332543
$ws.Range("C2").Value = -0.07804303399018395
$ws.Range("E2").Value = 0.1675937889528864

# Row 3
$ws.Range("B3").Value = 9.415194662350457
$ws.Range("C3").Value = 9.793741970092768
$ws.Range("E3").Value = 10.17228927783508

# Row 4
$ws.Range("B4").Value = -0.08861640112952046
$ws.Range("C4").Value = 0.4659205697519017
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 1.020457540633324

# Row 5
$ws.Range("B5").Value = -8.353335167575146
$ws.Range("C5").Value = -7.517289310840783
$ws.Range("E5").Value = -6.68124345410642

# Row 6
$ws.Range("B6").Value = 5.332219269048514
$ws.Range("C6").Value = 5.806335826943708
$ws.Range("E6").Value = 6.280452384838902

# Row 7
$ws.Range("B7").Value = -0.2542702753794868
$ws.Range("C7").Value = 0.06131248183902593
$ws.Range("E7").Value = 0.3768952390575387
